$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.3920720726908886
$ws.Range("J2").Value = 0.3920720726908886
$ws.Range("M2").Value = 26.05761233333333
$ws.Range("N2").Value = 78.172837
$ws.Range("O2").Value = 0.1375266232872619
$ws.Range("P2").Value = 0.1375266232872619
$ws.Range("Q2").Value = 29.74717046470544
$ws.Range("R2").Value = 267.724534182349
$ws.Range("S2").Value = 0.0539203482424158
$ws.Range("T2").Value = 0.0539203482424158

# Row 3
$ws.Range("I3").Value = 0.3920720726908886
$ws.Range("J3").Value = 0.3920720726908886
$ws.Range("O3").Value = 0.6296605108298998
$ws.Range("P3").Value = 0.6296605108298998
$ws.Range("S3").Value = 0.2468723015726825
$ws.Range("T3").Value = 0.2468723015726825

# Row 4
$ws.Range("I4").Value = 0.3920720726908886
$ws.Range("J4").Value = 0.3920720726908886
$ws.Range("M4").Value = 44.11180366666667
$ws.Range("N4").Value = 132.335411
$ws.Range("O4").Value = 0.2328128658828383
$ws.Range("P4").Value = 0.2328128658828383
$ws.Range("Q4").Value = 50.35769687537189
$ws.Range("R4").Value = 453.219271878347
$ws.Range("S4").Value = 0.09127942287579024
$ws.Range("T4").Value = 0.09127942287579024

# Row 5
$ws.Range("G5").Value = 1.770097666666667
$ws.Range("H5").Value = 5.310293000000001
$ws.Range("I5").Value = 0.6079279273091115
$ws.Range("J5").Value = 0.6079279273091115
$ws.Range("M5").Value = 26.05761233333333
$ws.Range("N5").Value = 78.172837
$ws.Range("O5").Value = 0.1375266232872619
$ws.Range("P5").Value = 0.1375266232872619
$ws.Range("Q5").Value = 46.1245187901379
$ws.Range("R5").Value = 415.1206691112411
$ws.Range("S5").Value = 0.08360627504484612
$ws.Range("T5").Value = 0.08360627504484612

# Row 6
$ws.Range("G6").Value = 1.770097666666667
$ws.Range("H6").Value = 5.310293000000001
$ws.Range("I6").Value = 0.6079279273091115
$ws.Range("J6").Value = 0.6079279273091115
$ws.Range("O6").Value = 0.6296605108298998
$ws.Range("P6").Value = 0.6296605108298998
$ws.Range("Q6").Value = 211.1793874449876
$ws.Range("R6").Value = 1900.614487004888
$ws.Range("S6").Value = 0.3827882092572173
$ws.Range("T6").Value = 0.3827882092572173

# Row 7
$ws.Range("G7").Value = 1.770097666666667
$ws.Range("H7").Value = 5.310293000000001
$ws.Range("I7").Value = 0.6079279273091115
$ws.Range("J7").Value = 0.6079279273091115
$ws.Range("M7").Value = 44.11180366666667
$ws.Range("N7").Value = 132.335411
$ws.Range("O7").Value = 0.2328128658828383
$ws.Range("P7").Value = 0.2328128658828383
$ws.Range("Q7").Value = 78.08220074282478
$ws.Range("R7").Value = 702.739806685423
$ws.Range("S7").Value = 0.141533443007048
$ws.Range("T7").Value = 0.141533443007048
